# Inserts a new weekly price record at row 92 of the "Cebollín" sheet,
# pushing the previously existing rows 92..219 down to 93..220.
#
# The new row 92 reuses the same market/category metadata (columns
# A,B,C,E,F,G,H,I,N,O,Q,R) as the record that used to sit at row 92, but
# carries a newer date and its own volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 92:219 down to 93:220, leaving row 92 free for the new record.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly record.
$ws.Range("A92").Value = 10
$ws.Range("B92").Value = "Vega Modelo de Temuco"
$ws.Range("C92").Value = "La Araucanía"
$ws.Range("D92").Value = 44482
$ws.Range("E92").Value = 9
$ws.Range("F92").Value = 100112037
$ws.Range("G92").Value = "Cebollín"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 40
$ws.Range("K92").Value = 7000
$ws.Range("L92").Value = 7000
$ws.Range("M92").Value = 7000
$ws.Range("N92").Value = "`$/docena de paquetes"
$ws.Range("O92").Value = "Provincia de Cautín"
$ws.Range("P92").Value = 583
$ws.Range("Q92").Value = 12
$ws.Range("R92").Value = "Hortaliza"
